$wb = $excel.ActiveWorkbook

$longMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/05915c2bcbbedc00fd858e3e83bf6ac1064c34d6/e2e/96d536ec-2b5c-4b5c-aef1-22b9f514769e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9be1b90724a32a5618c50f911e3252750d649ca8/e2e/96d536ec-2b5c-4b5c-aef1-22b9f514769e.md."
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9be1b90724a32a5618c50f911e3252750d649ca8/e2e/96d536ec-2b5c-4b5c-aef1-22b9f514769e.md"
$displayName = "96d536ec-2b5c-4b5c-aef1-22b9f514769e.md"

# ---- zh-cn sheet, row 7 (handback generated for the localized target) ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "96d536ec-2b5c-4b5c-aef1-22b9f514769e.3eef0f4811f4f1032be2197e8eebfd54c228ae3a.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-28 16:55:27"
$wsZh.Range("P7").Value = $longMessage

$wsZh.Range("I7").Value = $displayName
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $displayName)

# ---- de-de sheet, row 7 ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "96d536ec-2b5c-4b5c-aef1-22b9f514769e.3eef0f4811f4f1032be2197e8eebfd54c228ae3a.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-28 16:55:34"
$wsDe.Range("P7").Value = $longMessage

$wsDe.Range("I7").Value = $displayName
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $displayName)
